$wb = $excel.ActiveWorkbook

# Add the new shared string value on Tabelle3 (sheet3) at cell C4
$ws3 = $wb.Worksheets.Item("Tabelle3")
$ws3.Range("C4").Value = "Test Long String"
$ws3.Range("C5").Select()

# Make Tabelle3 the active sheet (tabSelected / activeTab)
$ws3.Activate()
